# Fruta / hortaliza, semanal
# Insert a new weekly record at row 361 (Macroferia Regional de Talca - Piña, Caramelo),
# shifting the existing rows 361-380 down to 362-381.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 361, pushing rows 361..380 down to 362..381
$ws.Rows.Item(361).Insert()

# Populate the newly inserted row 361 with the new weekly record
$ws.Cells.Item(361, 1).Value = 5
$ws.Cells.Item(361, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(361, 3).Value = "Maule"
$ws.Cells.Item(361, 4).Value = 45041
$ws.Cells.Item(361, 5).Value = 7
$ws.Cells.Item(361, 6).Value = "Fruta"
$ws.Cells.Item(361, 7).Value = 100108
$ws.Cells.Item(361, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(361, 9).Value = 100108005
$ws.Cells.Item(361, 10).Value = "Piña"
$ws.Cells.Item(361, 11).Value = "Caramelo"
$ws.Cells.Item(361, 12).Value = "Segunda"
$ws.Cells.Item(361, 13).Value = 250
$ws.Cells.Item(361, 14).Value = 19000
$ws.Cells.Item(361, 15).Value = 19000
$ws.Cells.Item(361, 16).Value = 19000
$ws.Cells.Item(361, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(361, 18).Value = "Ecuador"
$ws.Cells.Item(361, 19).Value = 1357
$ws.Cells.Item(361, 20).Value = 14
